# Sync attendance_reports: reorder the "Recorded By" (column G) names so
# that "System" is no longer forced to the front of the list - the first
# entry moves to the end (a right-rotation of the comma-separated names),
# except for rows whose list includes "admin@admin.com" which stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$changed = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }
    if (-not $val.StartsWith("System, ")) {
        continue
    }
    if ($val.Contains("admin@admin.com")) {
        continue
    }

    $parts = $val -split ", "
    $n = $parts.Count
    if ($n -lt 2) {
        continue
    }

    $last = $parts[$n - 1]
    $rest = $parts[0..($n - 2)]
    $newParts = @($last) + $rest
    $newVal = $newParts -join ", "

    $cell.Value2 = $newVal
    $changed++
}

Write-Output "Rows updated: $changed"
